$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to force pure-numeric-looking strings to be written as literal
# text (matching the source data, which stores prices/links as inline strings),
# by formatting it as Text and pasting values-only into the destination cells.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

$helper.Value = "1.001"
$helper.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$helper.Value = "243.61"
$helper.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$helper.Value = "0.6280"
$helper.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$helper.Value = "0.07465"
$helper.Copy()
$ws.Range("D8").PasteSpecial(-4163)
$helper.Value = "0.2921"
$helper.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$helper.Value = "23.07"
$helper.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$helper.Value = "0.07721"
$helper.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$helper.Value = "4.978"
$helper.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$helper.Value = "0.6681"
$helper.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$helper.Value = "82.61"
$helper.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$helper.Value = "0.000009333"
$helper.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$helper.Value = "6.027"
$helper.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$helper.Value = "222.83"
$helper.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$helper.Value = "7.120"
$helper.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$helper.Value = "1.003"
$helper.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$helper.Value = "160.32"
$helper.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$helper.Value = "0.1395"
$helper.Copy()
$ws.Range("D26").PasteSpecial(-4163)
$helper.Value = "8.495"
$helper.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$helper.Value = "17.91"
$helper.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$helper.Value = "1.503"
$helper.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$helper.Value = "0.05865"
$helper.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$helper.Value = "4.147"
$helper.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$helper.Value = "4.067"
$helper.Copy()
$ws.Range("D32").PasteSpecial(-4163)
$helper.Value = "0.7496"
$helper.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$helper.Value = "1.851"
$helper.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$helper.Value = "1.136"
$helper.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$helper.Value = "2.607"
$helper.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$helper.Value = "6.551"
$helper.Copy()
$ws.Range("D41").PasteSpecial(-4163)
$helper.Value = "0.8938"
$helper.Copy()
$ws.Range("D42").PasteSpecial(-4163)
$helper.Value = "1.003"
$helper.Copy()
$ws.Range("D43").PasteSpecial(-4163)
$helper.Value = "102.10"
$helper.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$helper.Value = "65.71"
$helper.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$helper.Value = "0.00000000123"
$helper.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$helper.Value = "0.07782"
$helper.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$helper.Value = "0.5092"
$helper.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$helper.Value = "0.4059"
$helper.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$helper.Value = "9.027"
$helper.Copy()
$ws.Range("D51").PasteSpecial(-4163)

# Helper no longer needed
$helper.Clear()

# Remaining cells are plain text (never mistaken for numbers by Excel) so a
# direct .Value assignment is sufficient and keeps their original style.
$ws.Range("D2").Value = "29.070.18"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.831.88"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").Value = "1.838.66"
$ws.Range("E12").Value = "  -0.14%  "
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  -0.18%  "
$ws.Range("E16").Value = "  -5.75%  "
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "29.075.59"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "2.074.44"
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("E20").Value = "  +1.99%  "
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +1.97%  "
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").Value = "  +12.72%  "
$ws.Range("E31").Value = "  +1.69%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("E34").Value = "  +1.42%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  -3.60%  "
$ws.Range("D38").Value = "1.229.11"
$ws.Range("E38").Value = "  -1.73%  "
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +3.13%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Value = "1.991.85"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E46").Value = "  +2.51%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("B48").Value = "XinFinNetwork"
$ws.Range("C48").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("E48").Value = "  +16.42%  "
$ws.Range("E49").Value = "  -0.58%  "
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  +1.82%  "
